# Barclaycard export format fix:
#  - "Name des Karteninhabers" header shortened to "Karteninhaber"
#  - "Händlerdetails" header shortened to "Details"
#  - the merchant-details data cell (previously its own "Händler" string)
#    now reuses the same "Details" text as its column header
#  - column B widened very slightly
#  - selection/scroll position updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L13").Value = "Karteninhaber"
$ws.Range("O13").Value = "Details"
$ws.Range("O14").Value = "Details"

$ws.Columns.Item(2).ColumnWidth = 26.6

$ws.Range("O12").Select()
